# Natmi following Dr Hou advice
#
# The sheet holds ligand/receptor pair scores between sending and target
# clusters for Tnfsf13 -> Sdc2. Previously only two combinations (ECs->ECs
# and ECs->FAPs... originally ECs->"sCs"-labelled column) were present.
# The update expands the table to the full 3x3 matrix of sending/target
# clusters (ECs, FAPs, sCs) with refreshed scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ligand = "Tnfsf13"
$receptor = "Sdc2"

# Each entry: sending cluster, target cluster, and the 16 numeric values
# for columns E..T (in column order).
$rows = @(
    @{ A = "ECs";  D = "ECs";  V = @(3, 1, 0.665239, 1.995717, 0.3721652871971906, 0.3721652871971906, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 0.8790769721013334, 7.911692748912, 0.003802718148054192, 0.003802718148054193) },
    @{ A = "ECs";  D = "FAPs"; V = @(3, 1, 0.665239, 1.995717, 0.3721652871971906, 0.3721652871971906, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 70.49087072276033, 634.417836504843, 0.3049299684518219, 0.3049299684518219) },
    @{ A = "ECs";  D = "sCs";  V = @(3, 1, 0.665239, 1.995717, 0.3721652871971906, 0.3721652871971906, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 14.66375794749167, 131.973821527425, 0.06343260059731448, 0.0634326005973145) },
    @{ A = "FAPs"; D = "ECs";  V = @(2, 0.6666666666666666, 0.3279213333333333, 0.983764, 0.1834542731230215, 0.1834542731230214, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 0.4333301156337778, 3.899971040704, 0.001874502855967246, 0.001874502855967246) },
    @{ A = "FAPs"; D = "FAPs"; V = @(2, 0.6666666666666666, 0.3279213333333333, 0.983764, 0.1834542731230215, 0.1834542731230214, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 34.74760246352844, 312.728422171756, 0.1503114547223069, 0.1503114547223069) },
    @{ A = "FAPs"; D = "sCs";  V = @(2, 0.6666666666666666, 0.3279213333333333, 0.983764, 0.1834542731230215, 0.1834542731230214, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 7.228318029788888, 65.0548622681, 0.03126831554474732, 0.03126831554474733) },
    @{ A = "sCs";  D = "ECs";  V = @(2, 0.6666666666666666, 0.7943223333333332, 2.382967, 0.444380439679788, 0.444380439679788, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 1.049653540545778, 9.446881864911999, 0.00454059962264903, 0.004540599622649031) },
    @{ A = "sCs";  D = "FAPs"; V = @(2, 0.6666666666666666, 0.7943223333333332, 2.382967, 0.444380439679788, 0.444380439679788, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 84.16895718862143, 757.5206146975929, 0.3640987435251253, 0.3640987435251254) },
    @{ A = "sCs";  D = "sCs";  V = @(2, 0.6666666666666666, 0.7943223333333332, 2.382967, 0.444380439679788, 0.444380439679788, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 17.50912142596388, 157.582092833675, 0.07574109653201366, 0.07574109653201368) }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $ligand
    $ws.Cells.Item($r, 3).Value = $receptor
    $ws.Cells.Item($r, 4).Value = $row.D

    $c = 5
    foreach ($v in $row.V) {
        $ws.Cells.Item($r, $c).Value = $v
        $c++
    }

    $r++
}
